$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2 currently holds one row: AFS / POSTED (row1).
# We need to insert three new rows above it (FMS/PARKED, FMS/POSTED, AFS/PARKED)
# pushing the existing row down to row 4.
$ws2.Rows.Item(1).Resize(3, 1).EntireRow.Insert()

$ws2.Range("A1").Value = "FMS"
$ws2.Range("B1").Value = "PARKED"
$ws2.Range("A2").Value = "FMS"
$ws2.Range("B2").Value = "POSTED"
$ws2.Range("A3").Value = "AFS"
$ws2.Range("B3").Value = "PARKED"

# Sheet1: keep header row + one data row (AFS / POSTED), remove the rest.
$ws1.Range("A2").Value = "AFS"
$ws1.Range("B2").Value = "POSTED"
$ws1.Range("A3:B4").EntireRow.Delete()

# Update the selections to match the target state (Sheet1 stays the active tab).
$ws2.Activate() | Out-Null
$ws2.Range("A4:B4").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("E8").Select() | Out-Null
